# Scheduled-runner refresh of market-price / leve-profit figures across
# several crafting-job sheets (currentAveragePrice* / LevePrice* /
# LeveProfit* columns, H:N). Values are plain numeric literals (no
# formulas in these sheets), so each touched cell is rewritten in place.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(5,8).Value = 2396
$ws_ALC.Cells.Item(5,9).Value = 2384.4443
$ws_ALC.Cells.Item(5,10).Value = 2500
$ws_ALC.Cells.Item(5,11).Value = 2384.4443
$ws_ALC.Cells.Item(5,12).Value = 2500
$ws_ALC.Cells.Item(5,13).Value = -2269.4443
$ws_ALC.Cells.Item(5,14).Value = -2730

$ws_ALC.Cells.Item(6,8).Value = 44764.777
$ws_ALC.Cells.Item(6,9).Value = 50260.5
$ws_ALC.Cells.Item(6,11).Value = 150781.5
$ws_ALC.Cells.Item(6,13).Value = -150669.5

$ws_ALC.Cells.Item(15,8).Value = 659836.8
$ws_ALC.Cells.Item(15,9).Value = 659836.8
$ws_ALC.Cells.Item(15,11).Value = 1979510.4
$ws_ALC.Cells.Item(15,13).Value = -1979341.4

$ws_ALC.Cells.Item(18,8).Value = 925.1539
$ws_ALC.Cells.Item(18,9).Value = 960.5833
$ws_ALC.Cells.Item(18,10).Value = 500
$ws_ALC.Cells.Item(18,11).Value = 960.5833
$ws_ALC.Cells.Item(18,12).Value = 500
$ws_ALC.Cells.Item(18,13).Value = -676.5833
$ws_ALC.Cells.Item(18,14).Value = -1068

$ws_ALC.Cells.Item(33,8).Value = 18736.393
$ws_ALC.Cells.Item(33,9).Value = 26679.422
$ws_ALC.Cells.Item(33,10).Value = 1967.7778
$ws_ALC.Cells.Item(33,11).Value = 26679.422
$ws_ALC.Cells.Item(33,12).Value = 1967.7778
$ws_ALC.Cells.Item(33,13).Value = -26450.422
$ws_ALC.Cells.Item(33,14).Value = -2425.7778

$ws_ALC.Cells.Item(62,8).Value = 2248.8462
$ws_ALC.Cells.Item(62,9).Value = 2227.4167
$ws_ALC.Cells.Item(62,11).Value = 2227.4167
$ws_ALC.Cells.Item(62,13).Value = -1603.4167

$ws_ALC.Cells.Item(65,8).Value = 2248.8462
$ws_ALC.Cells.Item(65,9).Value = 2227.4167
$ws_ALC.Cells.Item(65,11).Value = 11137.0835
$ws_ALC.Cells.Item(65,13).Value = -8017.083500000001

$ws_ALC.Cells.Item(137,8).Value = 5816.6914
$ws_ALC.Cells.Item(137,9).Value = 5469.7554
$ws_ALC.Cells.Item(137,10).Value = 6135.306
$ws_ALC.Cells.Item(137,11).Value = 16409.2662
$ws_ALC.Cells.Item(137,12).Value = 18405.918
$ws_ALC.Cells.Item(137,13).Value = -13859.2662
$ws_ALC.Cells.Item(137,14).Value = -23505.918

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(45,8).Value = 2289.923
$ws_ARM.Cells.Item(45,9).Value = 1419.2222
$ws_ARM.Cells.Item(45,10).Value = 4249
$ws_ARM.Cells.Item(45,11).Value = 1419.2222
$ws_ARM.Cells.Item(45,12).Value = 4249
$ws_ARM.Cells.Item(45,13).Value = -1042.2222
$ws_ARM.Cells.Item(45,14).Value = -5003

$ws_ARM.Cells.Item(63,8).Value = 3472.4285
$ws_ARM.Cells.Item(63,9).Value = 2901.5557
$ws_ARM.Cells.Item(63,11).Value = 2901.5557
$ws_ARM.Cells.Item(63,13).Value = -2215.5557

$ws_ARM.Cells.Item(66,8).Value = 3472.4285
$ws_ARM.Cells.Item(66,9).Value = 2901.5557
$ws_ARM.Cells.Item(66,11).Value = 14507.7785
$ws_ARM.Cells.Item(66,13).Value = -11075.7785

$ws_ARM.Cells.Item(102,8).Value = 1085.1818
$ws_ARM.Cells.Item(102,9).Value = 770.8889
$ws_ARM.Cells.Item(102,11).Value = 770.8889
$ws_ARM.Cells.Item(102,13).Value = 851.1111

$ws_ARM.Cells.Item(110,8).Value = 21740112
$ws_ARM.Cells.Item(110,9).Value = 1059.9
$ws_ARM.Cells.Item(110,10).Value = 166667140
$ws_ARM.Cells.Item(110,11).Value = 1059.9
$ws_ARM.Cells.Item(110,12).Value = 166667140
$ws_ARM.Cells.Item(110,13).Value = 985.0999999999999
$ws_ARM.Cells.Item(110,14).Value = -166671230

$ws_ARM.Cells.Item(122,8).Value = 3013.5625
$ws_ARM.Cells.Item(122,9).Value = 3270.75
$ws_ARM.Cells.Item(122,10).Value = 2242
$ws_ARM.Cells.Item(122,11).Value = 9812.25
$ws_ARM.Cells.Item(122,12).Value = 6726
$ws_ARM.Cells.Item(122,13).Value = -7362.25
$ws_ARM.Cells.Item(122,14).Value = -11626

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(22,8).Value = 566.4375
$ws_BSM.Cells.Item(22,9).Value = 601.6
$ws_BSM.Cells.Item(22,11).Value = 601.6
$ws_BSM.Cells.Item(22,13).Value = -428.6

$ws_BSM.Cells.Item(105,8).Value = 3168.5334
$ws_BSM.Cells.Item(105,9).Value = 2819.4
$ws_BSM.Cells.Item(105,10).Value = 3866.8
$ws_BSM.Cells.Item(105,11).Value = 2819.4
$ws_BSM.Cells.Item(105,12).Value = 3866.8
$ws_BSM.Cells.Item(105,13).Value = -1072.4
$ws_BSM.Cells.Item(105,14).Value = -7360.8

$ws_BSM.Cells.Item(134,8).Value = 2104.1614
$ws_BSM.Cells.Item(134,9).Value = 1348.8
$ws_BSM.Cells.Item(134,10).Value = 5251.5
$ws_BSM.Cells.Item(134,11).Value = 4046.4
$ws_BSM.Cells.Item(134,12).Value = 15754.5
$ws_BSM.Cells.Item(134,13).Value = -1511.4
$ws_BSM.Cells.Item(134,14).Value = -20824.5

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(31,8).Value = 13810.875
$ws_CRP.Cells.Item(31,9).Value = 15500
$ws_CRP.Cells.Item(31,10).Value = 13698.267
$ws_CRP.Cells.Item(31,11).Value = 15500
$ws_CRP.Cells.Item(31,12).Value = 13698.267
$ws_CRP.Cells.Item(31,13).Value = -15205
$ws_CRP.Cells.Item(31,14).Value = -14288.267

$ws_CRP.Cells.Item(34,8).Value = 13810.875
$ws_CRP.Cells.Item(34,9).Value = 15500
$ws_CRP.Cells.Item(34,10).Value = 13698.267
$ws_CRP.Cells.Item(34,11).Value = 15500
$ws_CRP.Cells.Item(34,12).Value = 13698.267
$ws_CRP.Cells.Item(34,13).Value = -15298
$ws_CRP.Cells.Item(34,14).Value = -14102.267

$ws_CRP.Cells.Item(99,8).Value = 1030551.9
$ws_CRP.Cells.Item(99,9).Value = 1356588.1
$ws_CRP.Cells.Item(99,10).Value = 2283.7693
$ws_CRP.Cells.Item(99,11).Value = 1356588.1
$ws_CRP.Cells.Item(99,12).Value = 2283.7693
$ws_CRP.Cells.Item(99,13).Value = -1355090.1
$ws_CRP.Cells.Item(99,14).Value = -5279.7693

$ws_CRP.Cells.Item(126,8).Value = 1030551.9
$ws_CRP.Cells.Item(126,9).Value = 1356588.1
$ws_CRP.Cells.Item(126,10).Value = 2283.7693
$ws_CRP.Cells.Item(126,11).Value = 4069764.3
$ws_CRP.Cells.Item(126,12).Value = 6851.3079
$ws_CRP.Cells.Item(126,13).Value = -4067294.3
$ws_CRP.Cells.Item(126,14).Value = -11791.3079

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(7,8).Value = 474.23077
$ws_CUL.Cells.Item(7,9).Value = 400
$ws_CUL.Cells.Item(7,10).Value = 480.41666
$ws_CUL.Cells.Item(7,11).Value = 1200
$ws_CUL.Cells.Item(7,12).Value = 1441.24998
$ws_CUL.Cells.Item(7,13).Value = -1088
$ws_CUL.Cells.Item(7,14).Value = -1665.24998

$ws_CUL.Cells.Item(68,8).Value = 4141.8535
$ws_CUL.Cells.Item(68,9).Value = 3408.5454
$ws_CUL.Cells.Item(68,10).Value = 4410.7334
$ws_CUL.Cells.Item(68,11).Value = 10225.6362
$ws_CUL.Cells.Item(68,12).Value = 13232.2002
$ws_CUL.Cells.Item(68,13).Value = -9414.636200000001
$ws_CUL.Cells.Item(68,14).Value = -14854.2002

$ws_CUL.Cells.Item(71,8).Value = 4141.8535
$ws_CUL.Cells.Item(71,9).Value = 3408.5454
$ws_CUL.Cells.Item(71,10).Value = 4410.7334
$ws_CUL.Cells.Item(71,11).Value = 30676.9086
$ws_CUL.Cells.Item(71,12).Value = 39696.6006
$ws_CUL.Cells.Item(71,13).Value = -26620.9086
$ws_CUL.Cells.Item(71,14).Value = -47808.6006

$ws_CUL.Cells.Item(80,8).Value = 3799.8
$ws_CUL.Cells.Item(80,10).Value = 4750
$ws_CUL.Cells.Item(80,12).Value = 14250
$ws_CUL.Cells.Item(80,14).Value = -16122

$ws_CUL.Cells.Item(83,8).Value = 3799.8
$ws_CUL.Cells.Item(83,10).Value = 4750
$ws_CUL.Cells.Item(83,12).Value = 42750
$ws_CUL.Cells.Item(83,14).Value = -52110

$ws_CUL.Cells.Item(92,8).Value = 75.40000000000001
$ws_CUL.Cells.Item(92,10).Value = 72.333336
$ws_CUL.Cells.Item(92,12).Value = 217.000008
$ws_CUL.Cells.Item(92,14).Value = -2713.000008

$ws_CUL.Cells.Item(97,8).Value = 810.5714
$ws_CUL.Cells.Item(97,10).Value = 358.33334
$ws_CUL.Cells.Item(97,12).Value = 1075.00002
$ws_CUL.Cells.Item(97,14).Value = -2067.00002

$ws_CUL.Cells.Item(129,8).Value = 8445029
$ws_CUL.Cells.Item(129,9).Value = 158964
$ws_CUL.Cells.Item(129,10).Value = 12906756
$ws_CUL.Cells.Item(129,11).Value = 476892
$ws_CUL.Cells.Item(129,12).Value = 38720268
$ws_CUL.Cells.Item(129,13).Value = -471892
$ws_CUL.Cells.Item(129,14).Value = -38730268

$ws_CUL.Cells.Item(131,8).Value = 827527
$ws_CUL.Cells.Item(131,9).Value = 904.2222
$ws_CUL.Cells.Item(131,11).Value = 2712.6666
$ws_CUL.Cells.Item(131,13).Value = 2327.3334

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(102,8).Value = 2490.1765
$ws_GSM.Cells.Item(102,9).Value = 1562.9286
$ws_GSM.Cells.Item(102,11).Value = 1562.9286
$ws_GSM.Cells.Item(102,13).Value = 59.07140000000004

$ws_GSM.Cells.Item(126,8).Value = 4933.3335
$ws_GSM.Cells.Item(126,9).Value = 3035.1667
$ws_GSM.Cells.Item(126,11).Value = 9105.500100000001
$ws_GSM.Cells.Item(126,13).Value = -6635.500100000001

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(16,8).Value = 4585.5713
$ws_LTW.Cells.Item(16,9).Value = 3683.1667
$ws_LTW.Cells.Item(16,11).Value = 3683.1667
$ws_LTW.Cells.Item(16,13).Value = -3513.1667

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 70/73: HQ price column (N) is retired in favour of a recomputed
# NQ profit (M) for this leve, so N is cleared and M is populated.
$ws_WVR.Cells.Item(70,8).Value = 20000
$ws_WVR.Cells.Item(70,9).Value = 20000
$ws_WVR.Cells.Item(70,10).Value = 0
$ws_WVR.Cells.Item(70,11).Value = 20000
$ws_WVR.Cells.Item(70,12).Value = 0
$ws_WVR.Cells.Item(70,14).ClearContents()
$ws_WVR.Cells.Item(70,13).Value = -19685

$ws_WVR.Cells.Item(73,8).Value = 20000
$ws_WVR.Cells.Item(73,9).Value = 20000
$ws_WVR.Cells.Item(73,10).Value = 0
$ws_WVR.Cells.Item(73,11).Value = 20000
$ws_WVR.Cells.Item(73,12).Value = 0
$ws_WVR.Cells.Item(73,14).ClearContents()
$ws_WVR.Cells.Item(73,13).Value = -18908

$ws_WVR.Cells.Item(81,8).Value = 1654.6
$ws_WVR.Cells.Item(81,9).Value = 1654.6
$ws_WVR.Cells.Item(81,11).Value = 3309.2
$ws_WVR.Cells.Item(81,13).Value = -2248.2

$ws_WVR.Cells.Item(84,8).Value = 1654.6
$ws_WVR.Cells.Item(84,9).Value = 1654.6
$ws_WVR.Cells.Item(84,11).Value = 16546
$ws_WVR.Cells.Item(84,13).Value = -11242

$ws_WVR.Cells.Item(100,8).Value = 652.5454999999999
$ws_WVR.Cells.Item(100,9).Value = 541.125
$ws_WVR.Cells.Item(100,11).Value = 1082.25
$ws_WVR.Cells.Item(100,13).Value = -541.25

$ws_WVR.Cells.Item(132,8).Value = 3713.7188
$ws_WVR.Cells.Item(132,9).Value = 2808.4285
$ws_WVR.Cells.Item(132,10).Value = 5442
$ws_WVR.Cells.Item(132,11).Value = 8425.2855
$ws_WVR.Cells.Item(132,12).Value = 16326
$ws_WVR.Cells.Item(132,13).Value = -5895.2855
$ws_WVR.Cells.Item(132,14).Value = -21386

$ws_WVR.Cells.Item(136,8).Value = 5043.0713
$ws_WVR.Cells.Item(136,9).Value = 3583.9707
$ws_WVR.Cells.Item(136,11).Value = 10751.9121
$ws_WVR.Cells.Item(136,13).Value = -8201.9121
